# Apply the vaccine-induced response term changes:
# - Cell A13 ("VO:0010461") is replaced with the new term "VO:0010463"
# - Update the view's scroll/selection state to match the saved workbook
#   (top-left cell reset to default, active cell/selection moved to A14)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the reserved-ID cell with the new term ID
$ws.Range("A13").Value = "VO:0010463"

# Reset scroll position (remove the frozen/scrolled topLeftCell offset) and
# move the active selection to A14
$appWindow = $excel.ActiveWindow
$appWindow.ScrollRow = 1
$appWindow.ScrollColumn = 1
$ws.Range("A14").Select()

$wb.Save()
